# Append three new cost-log rows (270-272) to the bottom of the data
# range on Sheet1, matching the "Github Auto Build at 2023-12-13 11:02"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A270").Value = "2023-12-13 11:01:55"
$ws.Range("B270").Value = 0.0004

$ws.Range("A271").Value = "2023-12-13 11:02:36"
$ws.Range("B271").Value = 0.0022

$ws.Range("A272").Value = "2023-12-13 11:02:54"
$ws.Range("B272").Value = 0.0006000000000000001
